# #12 Bubble boxes reshaped on solution slides on presentation
#
# Slide 18 ("Result") holds three narrative-rectangle callout bubbles
# (wedgeRectCallout) that get resized/repositioned and whose callout
# "pointer" adjustment guides (adj1/adj2) get tweaked.
#
# NOTE on the literal Left/Top/Width/Height numbers below: this COM host
# stores Shape.Left/Top/Width/Height as single-precision (32-bit) floats
# and truncates (rather than rounds) pt*12700 back to EMU on save. Using
# the "obvious" value (emu/12700) can therefore land 1 EMU short. The
# constants used here were solved so that, after the float32 round-trip
# and truncation this host performs, they reproduce the exact target EMU
# offsets/extents from the target OOXML.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(18)

# --- "Bulle narrative : rectangle 10" (box to show SAI understanding) ---
$shp1 = $s.Shapes.Item(9)
$shp1.Left = 212.6971673543307
$shp1.Top = 407.55952555905515
$shp1.Width = 144.66054918110237
$shp1.Height = 72.66054918110237
$shp1.Adjustments.Item(2) = -1.17696
$shp1.Adjustments.Item(1) = -0.45989

# --- "Bulle narrative : rectangle 11" (click to match the strategy...) ---
$shp2 = $s.Shapes.Item(10)
$shp2.Left = 529.34031296063
$shp2.Top = 409.1747284094488
$shp2.Width = 154.98992225984253
$shp2.Height = 65.03960829921259
$shp2.Adjustments.Item(2) = -2.47095
$shp2.Adjustments.Item(1) = -0.1898

# --- "Bulle narrative : rectangle 12" (show the understading) ---
$shp3 = $s.Shapes.Item(11)
$shp3.Left = 766.2385806771654
$shp3.Top = 404.0644074488189
$shp3.Width = 165.7981082362205
$shp3.Height = 51.71543407086614
$shp3.Adjustments.Item(2) = -1.59352
$shp3.Adjustments.Item(1) = -0.28992
